$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("D2").Value = "E6heiLs"
$ws.Range("D3").Value = "HCz67vr"
$ws.Range("D4").Value = "VQULTYp"
$ws.Range("D5").Value = "47dWcMW"
$ws.Range("D6").Value = "LO6VDaL"
$ws.Range("D7").Value = "tFIXHDy"
$ws.Range("D8").Value = "or5hkdJ"
$ws.Range("D9").Value = "rPT8s4K"
$ws.Range("D10").Value = "asQNTwb"
$ws.Range("D11").Value = "JPBuHoZ"
$ws.Range("D12").Value = "eF80d5y"
$ws.Range("D13").Value = "X7llUky"
$ws.Range("D14").Value = "vcOBslO"
$ws.Range("D15").Value = "z0hKsnc"
$ws.Range("D16").Value = "XgztlEH"
$ws.Range("D17").Value = "iT6SeRU"
$ws.Range("D18").Value = "nTscjpp"
$ws.Range("D19").Value = "wNecFKN"
$ws.Range("D20").Value = "t4evQI0"
$ws.Range("D21").Value = "upY07IW"
$ws.Range("D22").Value = "SXgPQsx"
$ws.Range("D23").Value = "n1mbLt9"
$ws.Range("D24").Value = "4kiiGdm"
$ws.Range("D25").Value = "op0YXGv"
$ws.Range("D26").Value = "HKGdLbK"
$ws.Range("D27").Value = "OVkbyQH"
$ws.Range("D28").Value = "rSiTxyV"
$ws.Range("D29").Value = "0FYUwts"
$ws.Range("D30").Value = "BogNRaD"
$ws.Range("D31").Value = "SWpJeTh"
